# Insert a new data row at row 932 (shifts existing rows 932:974 down to 933:975)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(932).Insert()

$ws.Cells.Item(932, 1).Value = 4
$ws.Cells.Item(932, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(932, 3).Value = "Los Lagos"
$ws.Cells.Item(932, 4).Value = 45147
$ws.Cells.Item(932, 5).Value = 10
$ws.Cells.Item(932, 6).Value = 100112004
$ws.Cells.Item(932, 7).Value = "Cebolla"
$ws.Cells.Item(932, 8).Value = "Sin especificar"
$ws.Cells.Item(932, 9).Value = "1a (guarda)"
$ws.Cells.Item(932, 10).Value = 150
$ws.Cells.Item(932, 11).Value = 15000
$ws.Cells.Item(932, 12).Value = 15000
$ws.Cells.Item(932, 13).Value = 15000
$ws.Cells.Item(932, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(932, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(932, 16).Value = 833
$ws.Cells.Item(932, 17).Value = 18
$ws.Cells.Item(932, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date number format used throughout column D
$ws.Cells.Item(932, 4).NumberFormat = $ws.Cells.Item(933, 4).NumberFormat
